$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "SW70"
$ws.Range("A4").Value = "SW100"

$ws.Range("B3").Value = 556987
$ws.Range("B4").Value = 996874

$ws.Range("C3").Value = "EG"
$ws.Range("C4").Value = "TUE"

$ws.Range("C5").Select()
